$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 177-194 were scraped without a "Quest Points" value, so the data
# in columns D (Series/N/A) and E (Release date) needs to shift right
# into E and F respectively, and D gets filled in with the missing quest
# points value of 0.
for ($r = 177; $r -le 194; $r++) {
    $series = $ws.Cells.Item($r, 4).Value2
    $releaseDate = $ws.Cells.Item($r, 5).Value2

    $ws.Cells.Item($r, 6).Value = $releaseDate
    $ws.Cells.Item($r, 5).Value = $series
    $ws.Cells.Item($r, 4).Value = 0
}
